$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with latest feed values.

$ws.Range('D2').Value = "'25.921.98"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "  +0.29%  "

$ws.Range('D3').Value = "'1.643.66"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "  +0.56%  "

$ws.Range('E4').Value = "  +0.13%  "

$ws.Range('D5').Value = "'215.11"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "  +0.06%  "

$ws.Range('D6').Value = "'0.5092"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "  +1.55%  "

$ws.Range('E7').Value = "  +0.13%  "

$ws.Range('D8').Value = "'0.2577"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = "  +0.34%  "

$ws.Range('E9').Value = "  +0.21%  "

$ws.Range('D10').Value = "'19.68"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "  +0.48%  "

$ws.Range('E11').Value = "  +0.76%  "

$ws.Range('D12').Value = "'4.311"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = "  +1.78%  "

$ws.Range('D13').Value = "'1.657.15"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "  +1.39%  "

$ws.Range('D14').Value = "'0.5469"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "  +0.82%  "

$ws.Range('D15').Value = "'0.0₅7884"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = "  -0.42%  "

$ws.Range('D16').Value = "'64.82"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "  +2.21%  "

$ws.Range('D17').Value = "'25.977.66"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "  +0.46%  "

$ws.Range('D18').Value = "'1.004"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "  +0.18%  "

$ws.Range('D19').Value = "'198.58"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "  -2.06%  "

$ws.Range('D20').Value = "'4.437"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = "  +2.78%  "

$ws.Range('E21').Value = "  +1.07%  "

$ws.Range('E22').Value = "  +1.44%  "

$ws.Range('D23').Value = "'1.005"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "  +0.21%  "

$ws.Range('D24').Value = "'1.854"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = "  -3.22%  "

$ws.Range('D25').Value = "'139.96"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "  -0.57%  "

$ws.Range('D26').Value = "'0.1147"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "  +0.57%  "

$ws.Range('D27').Value = "'6.892"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "  +2.92%  "

$ws.Range('D28').Value = "'15.76"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "  +0.56%  "

$ws.Range('D29').Value = "'1.237"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "  -0.03%  "

$ws.Range('D30').Value = "'0.05014"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "  +0.70%  "

$ws.Range('D31').Value = "'3.283"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "  +0.88%  "

$ws.Range('E32').Value = "  +0.89%  "

$ws.Range('D33').Value = "'1.541"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "  +0.31%  "

$ws.Range('D34').Value = "'2.360"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "  -0.14%  "

$ws.Range('D35').Value = "'0.8936"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "  +0.19%  "

$ws.Range('E36').Value = "  -1.24%  "

$ws.Range('D37').Value = "'1.134.21"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = "  -2.74%  "

$ws.Range('D38').Value = "'0.5533"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "  -1.35%  "

$ws.Range('D39').Value = "'0.01561"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "  +0.28%  "

$ws.Range('D40').Value = "'1.004"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "  +0.17%  "

$ws.Range('D41').Value = "'5.656"
$ws.Range('D41').Style = "Normal"

$ws.Range('D42').Value = "'0.8142"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "  +0.93%  "

$ws.Range('D43').Value = "'99.71"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "  +0.56%  "

$ws.Range('E44').Value = "  +8.32%  "

$ws.Range('D45').Value = "'1.784.30"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "  +0.71%  "

$ws.Range('D46').Value = "'0.4522"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "  +0.13%  "

$ws.Range('D47').Value = "'55.23"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = "  +0.97%  "

$ws.Range('D48').Value = "'1.005"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "  +0.25%  "

$ws.Range('D49').Value = "'0.05089"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "  +0.13%  "

$ws.Range('D50').Value = "'0.09572"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "  +3.50%  "

$ws.Range('D51').Value = "'1.003"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "  +0.04%  "
